# Realestate Update resale numbers 2023-06-21 21:41
# Appends the newest snapshot row (row 66) to the CityResaleNum sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CityResaleNum")

$newRow = 66

# Columns A-D hold text (Date/Time/Weekday/Week) and must remain plain text,
# not get auto-converted into Excel date/number serials.
$textValues = @{
    1 = "2023-06-21"  # A: Date
    2 = "21:33:02"    # B: Time
    3 = "Wednesday"   # C: Weekday
    4 = "25"          # D: Week
}

foreach ($col in $textValues.Keys) {
    $cell = $ws.Cells.Item($newRow, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $textValues[$col]
    $cell.ClearFormats()
}

# Columns E-T hold the numeric resale counts per city.
$numValues = @{
    5  = 122363  # E: Beijing
    6  = 133528  # F: Guangzhou
    7  = 162575  # G: Suzhou
    8  = 133622  # H: Hangzhou
    9  = 177314  # I: Nanjing
    10 = 114678  # J: Xi_an
    11 = 202269  # K: Chengdu
    12 = 225695  # L: Chongqing
    13 = 175436  # M: Tianjin
    14 = 103897  # N: Hefei
    15 = 39341   # O: Fuzhou
    16 = 33887   # P: Xiamen
    17 = 51958   # Q: Changsha
    18 = -1      # R: Shanghai
    19 = 36343   # S: Shenzhen
    20 = -1      # T: Wuhan
}

foreach ($col in $numValues.Keys) {
    $ws.Cells.Item($newRow, $col).Value = $numValues[$col]
}
